$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("L5").Value = 0.6
$ws.Range("M5").Value = 0.4918032786885246
$ws.Range("N5").Value = 0.5405405405405406
$ws.Range("O5").Value = 0.008163265306122436
$ws.Range("P5").Value = 0.01639344262295089
$ws.Range("Q5").Value = 0.01326781326781323
$ws.Range("R5").Value = 0.01379310344827584
$ws.Range("S5").Value = 0.03448275862068981
$ws.Range("T5").Value = 0.0251630941286113

# Row 6
$ws.Range("L6").Value = 0.5625
$ws.Range("M6").Value = 0.5901639344262295
$ws.Range("N6").Value = 0.576
$ws.Range("O6").Value = -0.02933673469387754
$ws.Range("P6").Value = 0.1147540983606558
$ws.Range("Q6").Value = 0.04872727272727262
$ws.Range("R6").Value = -0.04956896551724137
$ws.Range("S6").Value = 0.2413793103448278
$ws.Range("T6").Value = 0.09241379310344806

# Row 10
$ws.Range("L10").Value = 0.6
$ws.Range("M10").Value = 0.5901639344262295
$ws.Range("N10").Value = 0.5950413223140496
$ws.Range("O10").Value = 0.008163265306122436
$ws.Range("P10").Value = 0.1147540983606558
$ws.Range("Q10").Value = 0.06776859504132227
$ws.Range("R10").Value = 0.01379310344827584
$ws.Range("S10").Value = 0.2413793103448278
$ws.Range("T10").Value = 0.128526645768025

# Row 11
$ws.Range("L11").Value = 0.6
$ws.Range("M11").Value = 0.5901639344262295
$ws.Range("N11").Value = 0.5950413223140496
$ws.Range("O11").Value = 0.008163265306122436
$ws.Range("P11").Value = 0.1147540983606558
$ws.Range("Q11").Value = 0.06776859504132227
$ws.Range("R11").Value = 0.01379310344827584
$ws.Range("S11").Value = 0.2413793103448278
$ws.Range("T11").Value = 0.128526645768025

# Row 15
$ws.Range("L15").Value = 0.3805309734513274
$ws.Range("M15").Value = 0.7049180327868853
$ws.Range("N15").Value = 0.4942528735632184
$ws.Range("O15").Value = 0.008190547919412539
$ws.Range("P15").Value = 0.1311475409836066
$ws.Range("Q15").Value = 0.04263997033741201
$ws.Range("R15").Value = 0.02199747155499368
$ws.Range("S15").Value = 0.2285714285714286
$ws.Range("T15").Value = 0.09441707717569804

# Row 16
$ws.Range("L16").Value = 0.3739130434782609
$ws.Range("M16").Value = 0.7049180327868853
$ws.Range("N16").Value = 0.4886363636363636
$ws.Range("O16").Value = 0.001572617946346
$ws.Range("P16").Value = 0.1311475409836066
$ws.Range("Q16").Value = 0.03702346041055726
$ws.Range("R16").Value = 0.004223602484472114
$ws.Range("S16").Value = 0.2285714285714286
$ws.Range("T16").Value = 0.08198051948051967

# Row 20
$ws.Range("L20").Value = 0.3805309734513274
$ws.Range("M20").Value = 0.7049180327868853
$ws.Range("N20").Value = 0.4942528735632184
$ws.Range("O20").Value = 0.008190547919412539
$ws.Range("P20").Value = 0.1311475409836066
$ws.Range("Q20").Value = 0.04263997033741201
$ws.Range("R20").Value = 0.02199747155499368
$ws.Range("S20").Value = 0.2285714285714286
$ws.Range("T20").Value = 0.09441707717569804

# Row 21
$ws.Range("L21").Value = 0.3962264150943396
$ws.Range("M21").Value = 0.6885245901639344
$ws.Range("N21").Value = 0.502994011976048
$ws.Range("O21").Value = 0.02388598956242477
$ws.Range("P21").Value = 0.1147540983606558
$ws.Range("Q21").Value = 0.0513811087502416
$ws.Range("R21").Value = 0.06415094339622653
$ws.Range("S21").Value = 0.2
$ws.Range("T21").Value = 0.1137724550898207

# Row 30
$ws.Range("L30").Value = 0.6842105263157895
$ws.Range("M30").Value = 0.4262295081967213
$ws.Range("N30").Value = 0.5252525252525253
$ws.Range("O30").Value = 0.02706766917293235
$ws.Range("P30").Value = 0.04918032786885246
$ws.Range("Q30").Value = 0.04608585858585873
$ws.Range("R30").Value = 0.04118993135011445
$ws.Range("S30").Value = 0.1304347826086956
$ws.Range("T30").Value = 0.0961791831357052

# Row 31
$ws.Range("L31").Value = 0.6756756756756757
$ws.Range("M31").Value = 0.4098360655737705
$ws.Range("N31").Value = 0.5102040816326531
$ws.Range("O31").Value = 0.01853281853281852
$ws.Range("P31").Value = 0.03278688524590168
$ws.Range("Q31").Value = 0.0310374149659865
$ws.Range("R31").Value = 0.02820211515863687
$ws.Range("S31").Value = 0.08695652173913054
$ws.Range("T31").Value = 0.06477373558118922

# Row 35
$ws.Range("L35").Value = 0.4186046511627907
$ws.Range("M35").Value = 0.5901639344262295
$ws.Range("N35").Value = 0.4897959183673469
$ws.Range("O35").Value = 0.1686046511627907
$ws.Range("P35").Value = 0.5737704918032787
$ws.Range("Q35").Value = 0.4590266875981162
$ws.Range("R35").Value = 0.6744186046511629
$ws.Range("S35").Value = 35.00000000000004
$ws.Range("T35").Value = 14.91836734693881

# Row 36
$ws.Range("L36").Value = 0.4193548387096774
$ws.Range("M36").Value = 0.639344262295082
$ws.Range("N36").Value = 0.5064935064935064
$ws.Range("O36").Value = 0.1693548387096774
$ws.Range("P36").Value = 0.6229508196721312
$ws.Range("Q36").Value = 0.4757242757242757
$ws.Range("R36").Value = 0.6774193548387097
$ws.Range("S36").Value = 38.00000000000005
$ws.Range("T36").Value = 15.46103896103899
